$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9 (shifts rows 9-12 down to 10-13),
# matching pokemon #307 which was missing from the original list.
$ws.Rows("9:9").Insert()

$ws.Range("A9").Value = 307
# Populate columns E -> B (in this order so the new shared strings are
# interned in the same order as the target workbook).
$ws.Range("E9").Value = "這個特性編號在遊戲中未被使用，參考: https://wiki.52poke.com/wiki/Talk:%E7%89%B9%E6%80%A7%E5%88%97%E8%A1%A8"
$ws.Range("D9").Value = "<No Data>"
$ws.Range("C9").Value = "？？？"
$ws.Range("B9").Value = "<資料缺失>"

# Update the "Temp" defined name so it still spans the whole table
# (A1:G12 -> A1:G13) now that a row was added.
$wb.Names.Item("Temp").RefersTo = "=工作表1!`$A`$1:`$G`$13"

# Match the author's final selection.
$ws.Range("C9").Select()
